# Append a new daily row to the Gold Data sheet.
# Mirrors the source commit "Update Excel file from GitHub Actions":
# a new date ("23-11-2025") is recorded in column A, and since no new
# price paragraph was scraped that day, the price text from the most
# recent row is simply carried forward into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row (currently row 67) from the used range.
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$newDate = "23-11-2025"
$previousPrice = $ws.Cells.Item($lastRow, 2).Value2

$ws.Cells.Item($newRow, 1).Value = $newDate
$ws.Cells.Item($newRow, 2).Value = $previousPrice
